$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "grandes regiões e unidades da federação" sub-header in row 6 was a stray
# label with no data beneath it; every region row after it (norte, rondônia, ...)
# was actually holding the *next* region's numbers. Deleting row 6 outright
# shifts rows 7:38 up by one, so "norte" picks up the values that used to sit
# one row below it, all the way down to "distrito federal" taking on the old
# row 38 figures - and the sheet's used range shrinks from G38 to G37, with the
# now-unused shared string dropped automatically.
$ws.Rows("6:6").Delete()
